$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (TSM)
$ws.Range("D2").Value = 291.75
$ws.Range("E2").Value = 59
$ws.Range("F2").Value = 0.62
$ws.Range("K2").Value = 59.4
$ws.Range("N2").Value = 54.02451352198364

# Row 3 (ASML)
$ws.Range("D3").Value = 1112.59
$ws.Range("E3").Value = 63.2
$ws.Range("F3").Value = 6.88
$ws.Range("K3").Value = 58.2
$ws.Range("N3").Value = 54.02451352198364

# Row 4 (QCOM)
$ws.Range("D4").Value = 174.45
$ws.Range("E4").Value = 49.9
$ws.Range("F4").Value = 5.64
$ws.Range("H4").Value = 43
$ws.Range("I4").Value = 43
$ws.Range("K4").Value = 51.4
$ws.Range("N4").Value = 54.02451352198364

# Row 5 (NVDA)
$ws.Range("D5").Value = 183.56
$ws.Range("E5").Value = 46.7
$ws.Range("F5").Value = 1.83
$ws.Range("H5").Value = 70
$ws.Range("K5").Value = 49.2
$ws.Range("N5").Value = 54.02451352198364

# Row 6 (AMD)
$ws.Range("D6").Value = 215.87
$ws.Range("E6").Value = 31.5
$ws.Range("F6").Value = 0.76
$ws.Range("H6").Value = 60
$ws.Range("J6").Value = 66
$ws.Range("K6").Value = 46.2
$ws.Range("N6").Value = 54.02451352198364
